# Add a new "Hungary" worksheet (FC600/FC700 Hungary Market test data),
# modeled on the existing "Slovakia" sheet, and make it the active tab.

$wb = $excel.ActiveWorkbook

# The Slovakia sheet is the template for the new Hungary sheet.
$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate it, placing the copy immediately after Slovakia (i.e. at the end).
[void]$slovakia.Copy($null, $slovakia)

# The newly created copy is now the last sheet in the workbook.
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# Update the market name / ticket reference cells for Hungary.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3590/T3615"

# Slovakia is no longer the active sheet - reset its selection to a
# "whole sheet" selection like the other inactive sheets in the workbook.
[void]$slovakia.Range("A1:XFD1048576").Select()

# Hungary becomes the active sheet / active tab, selected at B10.
[void]$hungary.Activate()
[void]$hungary.Range("B10").Select()
